# Updates cryptos list figures (price + 1h volume change) to the new
# scrape snapshot, and fixes the Maker/EnergySwap row ordering (rows 43-44
# had swapped data in the prior snapshot).
#
# Price ("D") / Volume ("E") columns are stored as literal text in the
# workbook (not numbers), so any value that Excel's smart-entry would
# otherwise parse as a plain number is written with a leading apostrophe
# (quote-prefix) to force text, matching the original cell typing.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2: Bitcoin
$ws.Range("D2").Value = "41.687.49"

# Row 3: Ethereum
$ws.Range("D3").Value = "2.478.07"
$ws.Range("E3").Value = "  -0.12%  "

# Row 4: TetherUSD
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  +0.03%  "

# Row 5: BNB
$ws.Range("D5").Value = "'319.03"
$ws.Range("E5").Value = "  +1.66%  "

# Row 6: Solana
$ws.Range("D6").Value = "'93.08"
$ws.Range("E6").Value = "  -0.16%  "

# Row 7: XRP
$ws.Range("D7").Value = "'0.555"
$ws.Range("E7").Value = "  +2.14%  "

# Row 8: USDC
$ws.Range("E8").Value = "  -0.02%  "

# Row 9: Cardano
$ws.Range("E9").Value = "  +2.68%  "

# Row 10: Dogecoin
$ws.Range("D10").Value = "'0.0889"
$ws.Range("E10").Value = "  +13.33%  "

# Row 11: Avalanche
$ws.Range("E11").Value = "  +1.08%  "

# Row 12: TRON
$ws.Range("D12").Value = "'0.112"
$ws.Range("E12").Value = "  +0.90%  "

# Row 13: WrappedliquidstakedEther2.0
$ws.Range("D13").Value = "2.858.62"
$ws.Range("E13").Value = "  -0.26%  "

# Row 14: Polkadot
$ws.Range("D14").Value = "'6.96"

# Row 15: Chainlink
$ws.Range("D15").Value = "'15.72"
$ws.Range("E15").Value = "  -2.56%  "

# Row 16: WrappedEther
$ws.Range("D16").Value = "2.489.81"
$ws.Range("E16").Value = "  -1.03%  "

# Row 17: Polygon
$ws.Range("E17").Value = "  +4.11%  "

# Row 18: WrappedBTC
$ws.Range("D18").Value = "41.635.75"
$ws.Range("E18").Value = "  +0.24%  "

# Row 19: ShibaInu
$ws.Range("D19").Value = "0.0₃0963"
$ws.Range("E19").Value = "  +2.94%  "

# Row 21: Litecoin
$ws.Range("D21").Value = "'71.58"
$ws.Range("E21").Value = "  -0.13%  "

# Row 22: InternetComputer(DFINITY)
$ws.Range("D22").Value = "'11.55"
$ws.Range("E22").Value = "  +2.01%  "

# Row 23: BitcoinCash
$ws.Range("D23").Value = "'241.87"
$ws.Range("E23").Value = "  +2.26%  "

# Row 24: PancakeSwap
$ws.Range("E24").Value = "  +1.61%  "

# Row 25: ImmutableX
$ws.Range("E25").Value = "  +1.62%  "

# Row 26: Dai
$ws.Range("E26").Value = "  -0.01%  "

# Row 27: EthereumClassic
$ws.Range("D27").Value = "'24.92"
$ws.Range("E27").Value = "  +0.06%  "

# Row 28: Toncoin
$ws.Range("D28").Value = "'2.29"
$ws.Range("E28").Value = "  +4.12%  "

# Row 29: Cosmos
$ws.Range("E29").Value = "  +2.61%  "

# Row 30: InjectiveProtocol
$ws.Range("D30").Value = "'36.64"
$ws.Range("E30").Value = "  +2.27%  "

# Row 31: Monero
$ws.Range("D31").Value = "'156.85"
$ws.Range("E31").Value = "  -0.67%  "

# Row 32: Filecoin
$ws.Range("D32").Value = "'5.54"
$ws.Range("E32").Value = "  +1.02%  "

# Row 33: FirstDigitalUSD
$ws.Range("E33").Value = "  -0.12%  "

# Row 34: Hedera
$ws.Range("E34").Value = "  +2.49%  "

# Row 35: WEMIXToken
$ws.Range("E35").Value = "  -0.15%  "

# Row 36: Celestia
$ws.Range("D36").Value = "'17.56"
$ws.Range("E36").Value = "  +1.15%  "

# Row 37: LidoDAOToken
$ws.Range("E37").Value = "  +0.36%  "

# Row 38: ARBITRUM
$ws.Range("E38").Value = "  +0.47%  "

# Row 39: Stellar
$ws.Range("E39").Value = "  +1.43%  "

# Row 40: Kaspa
$ws.Range("E40").Value = "  -1.19%  "

# Row 41: RenderToken
$ws.Range("D41").Value = "'4.02"
$ws.Range("E41").Value = "  -1.90%  "

# Row 42: ApeXProtocol
$ws.Range("E42").Value = "  +2.21%  "

# Row 43: EnergySwap (was Maker row - data swapped with row 44)
$ws.Range("B43").Value = "EnergySwap"
$ws.Range("C43").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D43").Value = "'19.61"
$ws.Range("E43").Value = "  -0.85%  "

# Row 44: Maker (was EnergySwap row - data swapped with row 43)
$ws.Range("B44").Value = "Maker"
$ws.Range("C44").Value = "https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr"
$ws.Range("D44").Value = "1.984.51"
$ws.Range("E44").Value = "  +0.60%  "

# Row 45: VeChain
$ws.Range("D45").Value = "'0.0287"
$ws.Range("E45").Value = "  +0.88%  "

# Row 46: NEARProtocol
$ws.Range("D46").Value = "'3.04"
$ws.Range("E46").Value = "  +3.19%  "

# Row 47: FraxShare
$ws.Range("D47").Value = "'9.21"
$ws.Range("E47").Value = "  +0.93%  "

# Row 48: RocketPoolETH
$ws.Range("D48").Value = "2.713.15"
$ws.Range("E48").Value = "  -0.36%  "

# Row 49: Aave
$ws.Range("D49").Value = "'97.78"
$ws.Range("E49").Value = "  -0.13%  "

# Row 50: ordi
$ws.Range("D50").Value = "'68.16"
$ws.Range("E50").Value = "  +0.33%  "

# Row 51: BitcoinSV
$ws.Range("D51").Value = "'74.26"
$ws.Range("E51").Value = "  +2.73%  "
